$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate header row (row 2) from the "Product backlog" sheet/table.
# All subsequent rows shift up by one; the table range, autofilter, dimension and
# shared-string usage counts are recalculated automatically by Excel.
$ws.Rows("2:2").Delete()

# Update the active selection to match the post-edit state.
$ws.Range("A10").Select()
